$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 'orgtype'
$ws.Range("B22").Value = 'Individual'
$ws.Range("A23").Value = 'yearstart'
$ws.Range("B23").Value = '01-04-2023'
$ws.Range("A24").Value = 'yearend'
$ws.Range("B24").Value = '01-03-2024'
$ws.Range("A25").Value = 'orgstate'
$ws.Range("B25").Value = 'Karnataka'
$ws.Range("A26").Value = 'orgcity'
$ws.Range("B26").Value = 'Banglore'
$ws.Range("A27").Value = 'orgaddr'
$ws.Range("B27").Value = 'Business bay banglore near IT park'
$ws.Range("A28").Value = 'orgpincode'
$ws.Range("B28").Value = '411023'
$ws.Range("A29").Value = 'orgcountry'
$ws.Range("B29").Value = 'US'
$ws.Range("A30").Value = 'invflag'
$ws.Range("B30").Value = '1'
$ws.Range("A31").Value = 'invsflag'
$ws.Range("B31").Value = '19'
$ws.Range("A32").Value = 'billflag'
$ws.Range("B32").Value = 's'
$ws.Range("A33").Value = 'avflag'
$ws.Range("B33").Value = '2'
$ws.Range("A34").Value = 'invaliddcidname'
$ws.Range("B34").Value = '@#@@#,"   ",WQWE@#@#'
$ws.Range("A35").Value = 'invalidcustid'
$ws.Range("B35").Value = '757,574,274'
$ws.Range("A36").Value = 'invaliddcno'
$ws.Range("B36").Value = '102,110,232'
$ws.Range("A37").Value = 'invaliddcflag'
$ws.Range("B37").Value = '0,1,0'
$ws.Range("A38").Value = 'invalidtaxstate'
$ws.Range("B38").Value = 'maha,GJ,MP'

$ws.Range("C1").Select()
